$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the home/away match data (columns F:V) between a few row pairs ---
# Rows 68/69
$r68 = $ws.Range("F68:V68").Value()
$r69 = $ws.Range("F69:V69").Value()
$ws.Range("F68:V68").Value = $r69
$ws.Range("F69:V69").Value = $r68

# Rows 78/79
$r78 = $ws.Range("F78:V78").Value()
$r79 = $ws.Range("F79:V79").Value()
$ws.Range("F78:V78").Value = $r79
$ws.Range("F79:V79").Value = $r78

# Rows 172/173
$r172 = $ws.Range("F172:V172").Value()
$r173 = $ws.Range("F173:V173").Value()
$ws.Range("F172:V172").Value = $r173
$ws.Range("F173:V173").Value = $r172

# --- Append two new match rows (179, 180), copying formatting from row 178 ---
$ws.Range("A178:V178").Copy($ws.Range("A179:V180"))

$ws.Range("A179").Value = 178
$ws.Range("B179").Value = "spain"
$ws.Range("C179").Value = "laliga"
$ws.Range("D179").Value = "2023-2024"
$ws.Range("E179").Value = 45281.89583333334
$ws.Range("F179").Value = "Alaves"
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = "Real Madrid"
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = 4.44
$ws.Range("K179").Value = "10/12/2023 10:02"
$ws.Range("L179").Value = 6.04
$ws.Range("M179").Value = "21/12/2023 21:28"
$ws.Range("N179").Value = 3.85
$ws.Range("O179").Value = "10/12/2023 10:02"
$ws.Range("P179").Value = 4.13
$ws.Range("Q179").Value = "21/12/2023 21:27"
$ws.Range("R179").Value = 1.71
$ws.Range("S179").Value = "10/12/2023 10:02"
$ws.Range("T179").Value = 1.6
$ws.Range("U179").Value = "21/12/2023 21:19"
$ws.Range("V179").Value = "https://www.betexplorer.com/football/spain/laliga/alaves-real-madrid/bqUifoKa/"

$ws.Range("A180").Value = 179
$ws.Range("B180").Value = "spain"
$ws.Range("C180").Value = "laliga"
$ws.Range("D180").Value = "2023-2024"
$ws.Range("E180").Value = 45281.89583333334
$ws.Range("F180").Value = "Mallorca"
$ws.Range("G180").Value = 3
$ws.Range("H180").Value = "Osasuna"
$ws.Range("I180").Value = 2
$ws.Range("J180").Value = 2.22
$ws.Range("K180").Value = "10/12/2023 10:02"
$ws.Range("L180").Value = 2.26
$ws.Range("M180").Value = "21/12/2023 21:29"
$ws.Range("N180").Value = 3.02
$ws.Range("O180").Value = "10/12/2023 10:02"
$ws.Range("P180").Value = 2.92
$ws.Range("Q180").Value = "21/12/2023 21:29"
$ws.Range("R180").Value = 3.52
$ws.Range("S180").Value = "10/12/2023 10:02"
$ws.Range("T180").Value = 4.11
$ws.Range("U180").Value = "21/12/2023 21:29"
$ws.Range("V180").Value = "https://www.betexplorer.com/football/spain/laliga/mallorca-osasuna/CSRucmzs/"

Write-Host "Done"
